$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$summary = $wb.Worksheets.Item("Weekly Summary")

# --- Fill in new worklog rows 108-121 on the Logs sheet ---
$rows = @(
    @(108, 44439, 110, "Artefact Implementation", "Prepare detail report on experiment"),
    @(109, 44440, 100, "Artefact Implementation", "Re-design the codes for model save and reload to extend cross folding validation to 10"),
    @(110, 44441, 120, "Literature Search", "Search for detail of Real-time and preprocess data augmentation and select the most benefitical on for conducting experiment "),
    @(111, 44442, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, SDG, learing rate 0.1, 0.05, 0.001)"),
    @(112, 44442, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, SDG, learing rate 0.1, 0.05, 0.001)"),
    @(113, 44442, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, SDG, learing rate 0.1, 0.05, 0.001)"),
    @(114, 44442, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, SDG, learing rate 0.1, 0.05, 0.001)"),
    @(115, 44442, 90,  "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, SDG, learing rate 0.1, 0.05, 0.001)"),
    @(116, 44443, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, Adam, learing rate 0.01, 0.001, 0.0001)"),
    @(117, 44443, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, Adam, learing rate 0.01, 0.001, 0.0001)"),
    @(118, 44443, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, Adam, learing rate 0.01, 0.001, 0.0001)"),
    @(119, 44443, 120, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, Adam, learing rate 0.01, 0.001, 0.0001)"),
    @(120, 44443, 110, "Artefact Implementation", "Conduct experiment using 10 folds cross validation (CBAM + ResNet, Adam, learing rate 0.01, 0.001, 0.0001)"),
    @(121, 44444, 90,  "OnTrack Task",             "Task 7.1")
)

foreach ($row in $rows) {
    $r = $row[0]
    $logs.Range("A$r").Value = "Tithra Chap"
    $logs.Range("B$r").Value = "As. Prof. Richard Dazeley"
    $logs.Range("C$r").Value = "Emotion Recognition Using Facial Expression"
    $logs.Range("D$r").Value = $row[1]
    $logs.Range("E$r").Value = $row[2]
    $logs.Range("F$r").Value = $row[3]
    $logs.Range("G$r").Value = $row[4]
}

# --- Weekly Summary: Week9 totals (row 10) ---
$summary.Range("B10").Formula = "=ROUNDDOWN(SUM(Logs!E108:E121)/60,0)"
$summary.Range("C10").Formula = "=MOD(SUM(Logs!E108:E121),60)"

# --- View state updates ---
[void]$logs.Activate()
$excel.Goto($logs.Range("A107"), $false)
[void]$logs.Range("D114").Select()
[void]$summary.Activate()
[void]$summary.Range("C13").Select()
